$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "0000" number format to the whole B2:D28 block (creates style s="1"
# used by every cell in that range, including the ones left blank).
$ws.Range("B2:D28").NumberFormat = "0000"

# --- Column B (testbench) ---
$ws.Range("B2").Value  = "x"
$ws.Range("B3").Value  = "x"
$ws.Range("B4").Value  = "x"
$ws.Range("B5").Value  = "x"
$ws.Range("B6").Value  = "x"
$ws.Range("B9").Value  = "x"
$ws.Range("B10").Value = "x"
$ws.Range("B11").Value = "x"
$ws.Range("B12").Value = "x"
$ws.Range("B13").Value = "x"
$ws.Range("B14").Value = "x"
$ws.Range("B15").Value = "x"
$ws.Range("B16").Value = "/"
$ws.Range("B17").Value = "/"
$ws.Range("B18").Value = "x"
$ws.Range("B19").Value = "x"
$ws.Range("B20").Value = "x"
$ws.Range("B21").Value = "x"
$ws.Range("B22").Value = "x"
$ws.Range("B23").Value = "/"
$ws.Range("B24").Value = "x"
$ws.Range("B25").Value = "x"
$ws.Range("B26").Value = "x"
$ws.Range("B27").Value = "x"
$ws.Range("B28").Value = "x"

# --- Column C (program) ---
$ws.Range("C2").Value  = "x"
$ws.Range("C3").Value  = "x"
$ws.Range("C4").Value  = "x"
$ws.Range("C5").Value  = "x"
$ws.Range("C6").Value  = "x"
$ws.Range("C9").Value  = "x"
$ws.Range("C10").Value = "x"
$ws.Range("C11").Value = "x"
$ws.Range("C12").Value = "x"
$ws.Range("C13").Value = "x"
$ws.Range("C14").Value = "x"
$ws.Range("C15").Value = "x"
$ws.Range("C16").Value = "/"
$ws.Range("C17").Value = "/"
$ws.Range("C18").Value = "x"
$ws.Range("C19").Value = "x"
$ws.Range("C20").Value = "x"
$ws.Range("C21").Value = "x"
$ws.Range("C22").Value = "x"
$ws.Range("C23").Value = "/"
$ws.Range("C24").Value = "x"
$ws.Range("C25").Value = "x"
$ws.Range("C26").Value = "x"
$ws.Range("C27").Value = "x"
$ws.Range("C28").Value = "x"

# --- Column D (working) ---
$ws.Range("D2").Value  = "x"
$ws.Range("D3").Value  = "x"
$ws.Range("D4").Value  = "x"
$ws.Range("D5").Value  = "x"
$ws.Range("D6").Value  = "x"
$ws.Range("D16").Value = "x"
$ws.Range("D17").Value = "x"
$ws.Range("D18").Value = "x"
$ws.Range("D19").Value = "x"
$ws.Range("D20").Value = "x"
$ws.Range("D21").Value = "x"
$ws.Range("D22").Value = "x"
$ws.Range("D23").Value = "x"

# Move the active selection to D8 (was D24 before the edit).
$ws.Range("D8").Select()
